$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New relationship rows appended to the E:G mini-table (rows 30-32) ---
$ws.Range("E30").Value = "PrintCollection"
$ws.Range("F30").Value = "publishYear"
$ws.Range("G30").Value = "DateTime"

$ws.Range("E31").Value = "DigitalCollection"
$ws.Range("F31").Value = "productionYear"
$ws.Range("G31").Value = "DateTime"

$ws.Range("E32").Value = "Collection"
$ws.Range("F32").Value = "numOfCopies"
$ws.Range("G32").Value = "Integer"

# --- New rows appended to the A:C relationship table (rows 37-40) ---
$ws.Range("A37").Value = "Professor"
$ws.Range("B37").Value = "hasMentored"
$ws.Range("C37").Value = "AcademicPaper"

$ws.Range("A38").Value = "Author"
$ws.Range("B38").Value = "authorOf"
$ws.Range("C38").Value = "PrintCollection"

$ws.Range("A39").Value = "Director"
$ws.Range("B39").Value = "directorOf"
$ws.Range("C39").Value = "Movie"

$ws.Range("A40").Value = "Producer"
$ws.Range("B40").Value = "producerOf"
$ws.Range("C40").Value = "DigitalCollection"

# --- Match the author's final cursor position (selection moved to A41) ---
$ws.Range("A41").Select()
